$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the long step-by-step instruction texts in column D (Skenario
# Deskripsi) with short summary labels for each action (Tambah/View/Ubah/
# Hapus) against "Setup Parameter Investasi".
$ws.Range("D2").Value = "Tambah Setup Parameter Investasi"
$ws.Range("D3").Value = "View Setup Parameter Investasi"
$ws.Range("D4").Value = "Ubah Setup Parameter Investasi"
$ws.Range("D5").Value = "Hapus Setup Parameter Investasi"

# The rows no longer need to be as tall since the new text is much shorter;
# resize rows 2, 3 and 5 to 30pt and let row 4 auto-size back to the sheet
# default (its new text now fits on a single default-height line).
$ws.Rows(2).RowHeight = 30
$ws.Rows(3).RowHeight = 30
$ws.Rows(4).AutoFit()
$ws.Rows(5).RowHeight = 30

# Move the active selection to D5 (last edited cell).
$ws.Range("D5").Select() | Out-Null
